$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove trailing period from bullet-point descriptions in column E (rows 2-14)
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null -and $val.EndsWith(".")) {
        $cell.Value2 = $val.Substring(0, $val.Length - 1)
    }
}

$ws.Range("E20").Select()
